$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.913.94"
$ws.Range("E2").Value = "  +1.41%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.188.16"
$ws.Range("E3").Value = "  +1.05%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.48"
$ws.Range("E5").Value = "  +2.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.85"
$ws.Range("E6").Value = "  +2.72%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.186.63"
$ws.Range("E8").Value = "  +1.13%  "

$ws.Range("E9").Value = "  +1.23%  "

$ws.Range("E10").Value = "  +0.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.09"
$ws.Range("E11").Value = "  -0.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.512"
$ws.Range("E12").Value = "  +2.35%  "

$ws.Range("E13").Value = "  +2.37%  "

$ws.Range("E14").Value = "  +4.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.708.95"
$ws.Range("E15").Value = "  +0.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.898.08"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.40"
$ws.Range("E17").Value = "  +4.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.192.25"
$ws.Range("E18").Value = "  +0.47%  "

$ws.Range("E19").Value = "  +0.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "509.81"
$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.31"
$ws.Range("E21").Value = "  +3.44%  "

$ws.Range("E22").Value = "  +2.96%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.19"
$ws.Range("E23").Value = "  -0.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.97"
$ws.Range("E24").Value = "  +3.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.77"
$ws.Range("E25").Value = "  +0.43%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("E27").Value = "  +4.90%  "

$ws.Range("E28").Value = "  +2.47%  "

$ws.Range("E29").Value = "  +4.72%  "

$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.90"
$ws.Range("E30").Value = "  +10.68%  "

$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.88"
$ws.Range("E31").Value = "  +3.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.23"
$ws.Range("E32").Value = "  +2.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.23"
$ws.Range("E33").Value = "  +4.06%  "

$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.54"
$ws.Range("E35").Value = "  +0.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.76"
$ws.Range("E36").Value = "  -0.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0902"
$ws.Range("E37").Value = "  +0.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "480.57"
$ws.Range("E38").Value = "  +3.50%  "

$ws.Range("E39").Value = "  -2.16%  "

$ws.Range("E40").Value = "  -3.01%  "

$ws.Range("E41").Value = "  +1.74%  "

$ws.Range("E42").Value = "  +4.05%  "

$ws.Range("E43").Value = "  +4.96%  "

$ws.Range("E44").Value = "  +10.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.928.76"
$ws.Range("E45").Value = "  -4.52%  "

$ws.Range("E46").Value = "  -0.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.39"
$ws.Range("E47").Value = "  -0.86%  "

$ws.Range("E48").Value = "  -0.05%  "

$ws.Range("E49").Value = "  +1.67%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.33"
$ws.Range("E50").Value = "  +3.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.62"
$ws.Range("E51").Value = "  +6.69%  "
